$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '34.910.44'
$ws.Range("E2").Value = '  -0.64%  '
# Row 3
$ws.Range("D3").Value = '1.819.19'
$ws.Range("E3").Value = '  -0.76%  '
# Row 4
$ws.Range("E4").Value = '  -0.22%  '
# Row 5
Set-TextValue "D5" '230.53'
$ws.Range("E5").Value = '  -0.94%  '
# Row 6
Set-TextValue "D6" '0.617'
$ws.Range("E6").Value = '  +0.08%  '
# Row 7
$ws.Range("E7").Value = '  -0.28%  '
# Row 8
Set-TextValue "D8" '40.31'
$ws.Range("E8").Value = '  -6.12%  '
# Row 9
Set-TextValue "D9" '0.323'
$ws.Range("E9").Value = '  +4.32%  '
# Row 10
Set-TextValue "D10" '0.0683'
$ws.Range("E10").Value = '  -0.98%  '
# Row 11
Set-TextValue "D11" '0.0988'
$ws.Range("E11").Value = '  -1.81%  '
# Row 12
$ws.Range("D12").Value = '2.081.62'
$ws.Range("E12").Value = '  -0.85%  '
# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D13" '11.32'
$ws.Range("E13").Value = '  +0.95%  '
# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D14" '0.671'
$ws.Range("E14").Value = '  +0.89%  '
# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.810.83'
$ws.Range("E15").Value = '  -1.34%  '
# Row 16
Set-TextValue "D16" '4.62'
$ws.Range("E16").Value = '  -1.49%  '
# Row 17
$ws.Range("D17").Value = '34.881.98'
$ws.Range("E17").Value = '  -0.67%  '
# Row 18
Set-TextValue "D18" '69.71'
$ws.Range("E18").Value = '  -0.91%  '
# Row 19
$ws.Range("D19").Value = '0.0₃0785'
$ws.Range("E19").Value = '  -0.96%  '
# Row 20
Set-TextValue "D20" '240.99'
$ws.Range("E20").Value = '  +0.22%  '
# Row 21
Set-TextValue "D21" '12.06'
$ws.Range("E21").Value = '  +1.50%  '
# Row 22
Set-TextValue "D22" '4.70'
$ws.Range("E22").Value = '  +2.41%  '
# Row 23
$ws.Range("E23").Value = '  +0.05%  '
# Row 24
Set-TextValue "D24" '2.27'
$ws.Range("E24").Value = '  +1.27%  '
# Row 25
Set-TextValue "D25" '173.36'
$ws.Range("E25").Value = '  +1.07%  '
# Row 26
Set-TextValue "D26" '7.78'
$ws.Range("E26").Value = '  -0.60%  '
# Row 27
Set-TextValue "D27" '0.124'
$ws.Range("E27").Value = '  +2.09%  '
# Row 28
Set-TextValue "D28" '17.35'
$ws.Range("E28").Value = '  -1.19%  '
# Row 29
Set-TextValue "D29" '1.52'
$ws.Range("E29").Value = '  -4.36%  '
# Row 30
$ws.Range("E30").Value = '  -0.23%  '
# Row 31
Set-TextValue "D31" '4.01'
$ws.Range("E31").Value = '  +2.43%  '
# Row 32
Set-TextValue "D32" '0.0548'
$ws.Range("E32").Value = '  -1.20%  '
# Row 33
Set-TextValue "D33" '3.97'
$ws.Range("E33").Value = '  -1.10%  '
# Row 34
Set-TextValue "D34" '1.25'
$ws.Range("E34").Value = '  +12.40%  '
# Row 35
Set-TextValue "D35" '1.84'
$ws.Range("E35").Value = '  +1.85%  '
# Row 36
Set-TextValue "D36" '0.694'
$ws.Range("E36").Value = '  +2.00%  '
# Row 37
Set-TextValue "D37" '92.12'
$ws.Range("E37").Value = '  -1.64%  '
# Row 38
Set-TextValue "D38" '1.37'
$ws.Range("E38").Value = '  +7.33%  '
# Row 39
$ws.Range("D39").Value = '1.339.85'
$ws.Range("E39").Value = '  +0.97%  '
# Row 40
Set-TextValue "D40" '0.0193'
$ws.Range("E40").Value = '  -0.39%  '
# Row 41
Set-TextValue "D41" '0.977'
$ws.Range("E41").Value = '  -1.97%  '
# Row 42
$ws.Range("E42").Value = '  -2.98%  '
# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D43" '14.43'
$ws.Range("E43").Value = '  -3.42%  '
# Row 44
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D44" '2.41'
$ws.Range("E44").Value = '  -1.68%  '
# Row 45
$ws.Range("E45").Value = '  -1.56%  '
# Row 46
Set-TextValue "D46" '0.0520'
$ws.Range("E46").Value = '  +2.40%  '
# Row 47
Set-TextValue "D47" '6.20'
$ws.Range("E47").Value = '  -0.67%  '
# Row 48
$ws.Range("D48").Value = '1.994.94'
$ws.Range("E48").Value = '  -0.67%  '
# Row 49
$ws.Range("E49").Value = '  -0.32%  '
# Row 50
Set-TextValue "D50" '0.0663'
$ws.Range("E50").Value = '  +3.53%  '
# Row 51
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D51" '96.91'
$ws.Range("E51").Value = '  -4.07%  '
